$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.374.99"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "1.858.54"
$ws.Range("E3").Value = "  -3.76%  "
$ws.Range("E4").Value = "  -1.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4534"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3874"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07937"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.54%  "
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("D13").Value = "1.860.06"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.926"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.135"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.20%  "
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06519"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.542"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.17%  "
$ws.Range("D23").Value = "27.378.28"
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "2.090.38"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.082"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.439"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.488"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09297"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9379"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.598"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.269"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02240"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.231"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.197"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5923"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1902"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.281"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5621"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.378"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.926"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -1.41%  "
